$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestDataMappingSheet_SD")

# Insert a new row at row 44, shifting rows 44:104 down by one.
$ws.Rows.Item(44).Insert(-4121)  # xlShiftDown = -4121

# Copy formatting from the row above (row 43) so the new row matches
# the surrounding "FolioXXX" group styling (style index 4 / customFormat).
$ws.Rows.Item(43).Copy()
$ws.Rows.Item(44).PasteSpecial(-4122)  # xlPasteFormats = -4122

# Populate the new row's values.
$ws.Range("A44").Value = "FolioApprovalAndAuditHistory"
$ws.Range("B44").Value = "cares\Folio.xlsx"
$ws.Range("C44").Value = "FolioApprovalandAuditHistory"
$ws.Range("D44").Value = 1
